# OpenHornet Interconnect BOM update
# - Added cable for 4A5A2D1 APU lamp to interconnect:
#     - Added 4A5A2W2P2 (SMP-02V-BC) to 4A5A2W2 cable.
#     - Added 4A5A2D1P1 (SMR-02V-B) to 4A5A2D1 APU LAMP.
# - Updated connector quantities in Interconnect BOM:
#     - SMR-02V-B        : 7   -> 8
#     - SYM-001T-P0.6(N) : 125 -> 127
#     - SMP-02V-BC       : 7   -> 8
#     - SHF-001T-0.8BS   : 134 -> 136

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76: FIND NO 70 - SHF-001T-0.8BS - QTY 134 -> 136
$ws.Range("B76").Value = 136

# Row 77: FIND NO 71 - SMP-02V-BC - QTY 7 -> 8, add 4A5A2W2P2 to reference designators
$ws.Range("B77").Value = 8
$ws.Range("E77").Value = "3A2A1W1J3, 3A2A1W2J1, 3A2A1W2J2, 4A4A2W1J1, 4A4A2W1J2, 4A5A2W2P2, 4A7A2W1J2, 5A2W1P7"

# Row 82: FIND NO 76 - SMR-02V-B - QTY 7 -> 8, add 4A5A2D1P1 to reference designators
$ws.Range("B82").Value = 8
$ws.Range("E82").Value = "3A2A1W4P1, 3A2A1W5P1, 3A2A1W6P1, 4A4A3W1P1, 4A5A2D1P1, 4A8A2W1P1, 4A8A3W1P1, 5A2A3W1J1"

# Row 88: FIND NO 82 - SYM-001T-P0.6(N) - QTY 125 -> 127, add 4A5A2D1P1 to reference designators
$ws.Range("B88").Value = 127
$ws.Range("E88").Value = "2A2A1A3W1J1, 2A2A1A3W1J?, 2A5A1W3J1, 2A10W1P1, 2A13W3P1, 2A15W2P2, 3A2A1W3P1, 3A2A1W4P1, 3A2A1W5P1, 3A2A1W6P1, 4A2A1A2W1P1, 4A3A2W1P1, 4A4A1W1P1, 4A4A3W1P1, 4A5A2W2P1, 4A5A2D1P1, 4A7A3W1P1, 4A7A4J1, 4A8A1W1P1, 4A8A2W1P1, 4A8A3W1P1, 4A10W1P1, 5A2A1W1J1, 5A2A2W1J1, 5A2A3W1J1, 5A5A2W1J1, 5A10A1W1J1, 5A10A2W1J1, 10A1A1W1P1, 10A1A1W2P1"

# Move the active selection to where the edits were made (C76), matching the saved file state
$ws.Range("C76").Select() | Out-Null
